# ajustando filtro de datas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix rows 77-80: column F ("Data") was stored as text "15/09/2025"; convert
# to a real date serial (matching the date format already used by the other
# rows in the sheet).
for ($r = 77; $r -le 80; $r++) {
    $ws.Cells.Item($r, 6).Value2 = 45915
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# Append new consumption records (rows 81-84) for "admin" in "Consultórios"
$newRows = @(
    @{ A = "admin"; B = "seringa";  C = 23;  D = "Consultórios"; E = "15/09/2025"; F = "11:24:43" },
    @{ A = "admin"; B = "algodão";  C = 333; D = "Consultórios"; E = "15/09/2025"; F = "11:24:43" },
    @{ A = "admin"; B = "gazes";    C = 44;  D = "Consultórios"; E = "15/09/2025"; F = "11:24:43" },
    @{ A = "admin"; B = "luvas";    C = 55;  D = "Consultórios"; E = "15/09/2025"; F = "11:24:43" }
)

$row = 81
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 1).Value2 = $rec.A
    $ws.Cells.Item($row, 2).Value2 = $rec.B
    $ws.Cells.Item($row, 3).Value2 = $rec.C
    $ws.Cells.Item($row, 4).Value2 = $rec.D
    $ws.Cells.Item($row, 5).Value2 = $rec.E
    $ws.Cells.Item($row, 6).Value2 = $rec.F
    $row++
}
